$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.43270897865295
$ws.Range("C2").Value = 23.54311895370483

$ws.Range("B3").Value = 18.78989720344543
$ws.Range("C3").Value = 23.19322800636292

$ws.Range("B4").Value = 16.78513097763062
$ws.Range("C4").Value = 40.35292816162109

$ws.Range("B5").Value = 17.47561025619507
$ws.Range("C5").Value = 40.89450240135193

$ws.Range("B6").Value = 15.91819500923157
$ws.Range("C6").Value = 23.95486831665039

$ws.Range("B7").Value = 15.52881503105164
$ws.Range("C7").Value = 13.12521290779114

$ws.Range("B8").Value = 15.44611811637878
$ws.Range("C8").Value = 13.20347595214844

$ws.Range("B9").Value = 16.37338423728943
$ws.Range("C9").Value = 15.67780590057373

$ws.Range("B10").Value = 17.06266522407532
$ws.Range("C10").Value = 14.08163928985596

$ws.Range("B11").Value = 18.0543212890625
$ws.Range("C11").Value = 14.3774573802948

$ws.Range("B12").Value = 16.45379877090454
$ws.Range("C12").Value = 23.45046162605286

$ws.Range("B13").Value = 15.96139883995056
$ws.Range("C13").Value = 13.53701448440552

$ws.Range("B14").Value = 16.55931401252747
$ws.Range("C14").Value = 22.6620945930481
